$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 276.57144
$ws.Range("I2").Value = 276.57144
$ws.Range("K2").Value = 276.57144
$ws.Range("M2").Value = -163.57144

$ws.Range("H21").Value = 4016
$ws.Range("I21").Value = 4016
$ws.Range("K21").Value = 4016
$ws.Range("M21").Value = -3548

$ws.Range("H23").Value = 4016
$ws.Range("I23").Value = 4016
$ws.Range("K23").Value = 4016
$ws.Range("M23").Value = -3782

$ws.Range("H33").Value = 404.5
$ws.Range("I33").Value = 406.9091
$ws.Range("J33").Value = 400.7143
$ws.Range("K33").Value = 406.9091
$ws.Range("L33").Value = 400.7143
$ws.Range("M33").Value = -177.9091
$ws.Range("N33").Value = -858.7143

$ws.Range("N62").ClearContents()
$ws.Range("H62").Value = 4999
$ws.Range("I62").Value = 4999
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4999
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -4375

$ws.Range("N65").ClearContents()
$ws.Range("H65").Value = 4999
$ws.Range("I65").Value = 4999
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 24995
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -21875

$ws.Range("H107").Value = 413.22223
$ws.Range("I107").Value = 474.2857
$ws.Range("J107").Value = 199.5
$ws.Range("K107").Value = 474.2857
$ws.Range("L107").Value = 199.5
$ws.Range("M107").Value = 1445.7143
$ws.Range("N107").Value = -4039.5

$ws.Range("H137").Value = 2227
$ws.Range("I137").Value = 1558.4517
$ws.Range("J137").Value = 3707.3572
$ws.Range("K137").Value = 4675.355100000001
$ws.Range("L137").Value = 11122.0716
$ws.Range("M137").Value = -2125.355100000001
$ws.Range("N137").Value = -16222.0716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 21249.166
$ws.Range("I3").Value = 7835
$ws.Range("J3").Value = 34663.332
$ws.Range("K3").Value = 7835
$ws.Range("L3").Value = 34663.332
$ws.Range("M3").Value = -7720
$ws.Range("N3").Value = -34893.332

$ws.Range("H32").Value = 7129.58
$ws.Range("I32").Value = 4425.507
$ws.Range("J32").Value = 15241.8
$ws.Range("K32").Value = 4425.507
$ws.Range("L32").Value = 15241.8
$ws.Range("M32").Value = -4138.507
$ws.Range("N32").Value = -15815.8

$ws.Range("H74").Value = 5282.3447
$ws.Range("I74").Value = 5492.148
$ws.Range("K74").Value = 5492.148
$ws.Range("M74").Value = -4618.148

$ws.Range("H77").Value = 5282.3447
$ws.Range("I77").Value = 5492.148
$ws.Range("K77").Value = 27460.74
$ws.Range("M77").Value = -23092.74

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2977947.5
$ws.Range("J105").Value = 2647
$ws.Range("L105").Value = 2647
$ws.Range("N105").Value = -6141

$ws.Range("H122").Value = 4366.3228
$ws.Range("I122").Value = 4336.26
$ws.Range("K122").Value = 13008.78
$ws.Range("M122").Value = -10558.78

$ws.Range("H141").Value = 225178
$ws.Range("J141").Value = 320211.47
$ws.Range("L141").Value = 320211.47
$ws.Range("N141").Value = -330571.47

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10599.8
$ws.Range("J3").Value = 16666.334
$ws.Range("L3").Value = 49999.00199999999
$ws.Range("N3").Value = -50223.00199999999

$ws.Range("H5").Value = 64724.25
$ws.Range("I5").Value = 84239.914
$ws.Range("J5").Value = 6177.25
$ws.Range("K5").Value = 252719.742
$ws.Range("L5").Value = 18531.75
$ws.Range("M5").Value = -252607.742
$ws.Range("N5").Value = -18755.75

$ws.Range("H39").Value = 2245.4546
$ws.Range("J39").Value = 8000
$ws.Range("L39").Value = 24000
$ws.Range("N39").Value = -24588

$ws.Range("H40").Value = 82.92856999999999
$ws.Range("I40").Value = 82.92856999999999
$ws.Range("K40").Value = 331.71428
$ws.Range("M40").Value = -262.71428

$ws.Range("H41").Value = 2000
$ws.Range("I41").Value = 2000
$ws.Range("K41").Value = 6000
$ws.Range("M41").Value = -5662

$ws.Range("H135").Value = 64724.25
$ws.Range("I135").Value = 84239.914
$ws.Range("J135").Value = 6177.25
$ws.Range("K135").Value = 758159.226
$ws.Range("L135").Value = 55595.25
$ws.Range("M135").Value = -755624.226
$ws.Range("N135").Value = -60665.25

$ws.Range("H139").Value = 5998
$ws.Range("I139").Value = 1397.4
$ws.Range("K139").Value = 4192.200000000001
$ws.Range("M139").Value = 947.7999999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 123.125
$ws.Range("I2").Value = 83.25
$ws.Range("J2").Value = 163
$ws.Range("K2").Value = 83.25
$ws.Range("L2").Value = 163
$ws.Range("M2").Value = 29.75
$ws.Range("N2").Value = -389

$ws.Range("H31").Value = 7086.5
$ws.Range("I31").Value = 3099
$ws.Range("J31").Value = 34999
$ws.Range("K31").Value = 3099
$ws.Range("L31").Value = 34999
$ws.Range("M31").Value = -2807
$ws.Range("N31").Value = -35583

$ws.Range("H37").Value = 7086.5
$ws.Range("I37").Value = 3099
$ws.Range("J37").Value = 34999
$ws.Range("K37").Value = 3099
$ws.Range("L37").Value = 34999
$ws.Range("M37").Value = -2822
$ws.Range("N37").Value = -35553

$ws.Range("H122").Value = 55672.824
$ws.Range("J122").Value = 5123.25
$ws.Range("L122").Value = 15369.75
$ws.Range("N122").Value = -20269.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1892.6666
$ws.Range("I22").Value = 1965.7778
$ws.Range("J22").Value = 1783
$ws.Range("K22").Value = 1965.7778
$ws.Range("L22").Value = 1783
$ws.Range("M22").Value = -1670.7778
$ws.Range("N22").Value = -2373

$ws.Range("H27").Value = 1892.6666
$ws.Range("I27").Value = 1965.7778
$ws.Range("J27").Value = 1783
$ws.Range("K27").Value = 1965.7778
$ws.Range("L27").Value = 1783
$ws.Range("M27").Value = -1858.7778
$ws.Range("N27").Value = -1997

$ws.Range("H46").Value = 1274.5
$ws.Range("I46").Value = 750
$ws.Range("J46").Value = 1449.3334
$ws.Range("K46").Value = 750
$ws.Range("L46").Value = 1449.3334
$ws.Range("M46").Value = -562
$ws.Range("N46").Value = -1825.3334

$ws.Range("H136").Value = 3003
$ws.Range("I136").Value = 3002.3333
$ws.Range("K136").Value = 9006.999899999999
$ws.Range("M136").Value = -6456.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 2560.5
$ws.Range("I74").Value = 2639
$ws.Range("J74").Value = 2429.6667
$ws.Range("K74").Value = 2639
$ws.Range("L74").Value = 2429.6667
$ws.Range("M74").Value = -1703
$ws.Range("N74").Value = -4301.6667

$ws.Range("H77").Value = 2560.5
$ws.Range("I77").Value = 2639
$ws.Range("J77").Value = 2429.6667
$ws.Range("K77").Value = 7917
$ws.Range("L77").Value = 7289.000100000001
$ws.Range("M77").Value = -3237
$ws.Range("N77").Value = -16649.0001

$ws.Range("H136").Value = 1564.4193
$ws.Range("I136").Value = 1327.4828
$ws.Range("K136").Value = 3982.4484
$ws.Range("M136").Value = -1432.4484
